$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# Header text: bump the report's volume/issue number and the covered week
# -----------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# -----------------------------------------------------------------------
# Helper functions used below to move values between the "no data yet"
# placeholder cells (shared text "0" / "***.*") and real numeric cells,
# while keeping each cell's number format/font/alignment correct by
# borrowing it (via PasteSpecial Formats) from a same-column donor cell
# that already carries the desired style.
# -----------------------------------------------------------------------
function Set-PlaceholderText($addr, $text, $refAddr) {
    # Build the target cell as a literal text value referencing the shared string
    # table the same way the existing placeholder cells do, then restore the
    # donor cells number format/font/alignment (style) without touching the value.
    $dst = $ws.Range($addr)
    $dst.Formula = "=""" + $text + """"
    $dst.Copy()
    $dst.PasteSpecial(-4163)   # xlPasteValues: collapse formula to literal text
    $ref = $ws.Range($refAddr)
    $ref.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats: adopt the donors style
    $ws.Application.CutCopyMode = $false
}

function Set-NumFromPlaceholder($addr, $val, $refAddr) {
    $dst = $ws.Range($addr)
    $ref = $ws.Range($refAddr)
    $ref.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats: adopt the donors numeric style
    $ws.Application.CutCopyMode = $false
    $dst.Value = $val
}

Set-PlaceholderText "G14" "0" "C14"
Set-PlaceholderText "H14" "***.*" "L14"
Set-NumFromPlaceholder "D15" 2 "G15"
Set-NumFromPlaceholder "E15" -100 "H15"
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -80
$ws.Range("L15").Value = -90.909090909090
$ws.Range("M15").Value = -66.666666666666
Set-NumFromPlaceholder "C16" 6 "D16"
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -17.647058823529
$ws.Range("I16").Value = 52
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 15.555555555555
$ws.Range("L16").Value = 4
$ws.Range("M16").Value = 23.809523809523
$ws.Range("N16").Value = -73.869346733668
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -42.857142857142
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -4.761904761904
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 61
$ws.Range("K17").Value = -9.836065573770
$ws.Range("L17").Value = -5.172413793103
$ws.Range("M17").Value = 57.142857142857
$ws.Range("N17").Value = -20.289855072463
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -77.777777777777
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -57.142857142857
$ws.Range("I18").Value = 37
$ws.Range("J18").Value = 57
$ws.Range("K18").Value = -35.087719298245
$ws.Range("L18").Value = -21.276595744680
$ws.Range("M18").Value = -13.953488372093
$ws.Range("N18").Value = -89.181286549707
$ws.Range("C19").Value = 19
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 55
$ws.Range("G19").Value = 64
$ws.Range("H19").Value = -14.0625
$ws.Range("I19").Value = 233
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 16.5
$ws.Range("L19").Value = 23.936170212766
$ws.Range("M19").Value = 150.537634408602
$ws.Range("N19").Value = 89.430894308943
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 24
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 113
$ws.Range("J20").Value = 102
$ws.Range("K20").Value = 10.784313725490
$ws.Range("L20").Value = 22.826086956521
$ws.Range("M20").Value = 213.888888888889
$ws.Range("N20").Value = -77.972709551656
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = -23.404255319148
$ws.Range("F21").Value = 118
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = -20.805369127516
$ws.Range("I21").Value = 492
$ws.Range("J21").Value = 472
$ws.Range("K21").Value = 4.237288135593
$ws.Range("L21").Value = 10.313901345291
$ws.Range("M21").Value = 92.941176470588
$ws.Range("N21").Value = -60.828025477707
Set-PlaceholderText "C22" "0" "N22"
Set-NumFromPlaceholder "D22" 2 "G22"
Set-NumFromPlaceholder "E22" -100 "H22"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = -75
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -83.333333333333
Set-NumFromPlaceholder "C23" 3 "F23"
Set-NumFromPlaceholder "D23" 4 "G23"
Set-NumFromPlaceholder "E23" -25 "H23"
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 24
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -4
$ws.Range("M23").Value = 60
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = -4.301075268817
$ws.Range("I24").Value = 334
$ws.Range("J24").Value = 318
$ws.Range("K24").Value = 5.031446540880
$ws.Range("L24").Value = 3.086419753086
$ws.Range("M24").Value = 36.885245901639
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -27.272727272727
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 27.777777777777
$ws.Range("I25").Value = 140
$ws.Range("J25").Value = 127
$ws.Range("K25").Value = 10.236220472440
$ws.Range("L25").Value = -18.128654970760
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -60
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = -9.677419354838
$ws.Range("I26").Value = 98
$ws.Range("J26").Value = 92
$ws.Range("K26").Value = 6.521739130434
$ws.Range("L26").Value = 10.112359550561
$ws.Range("M26").Value = -2
Set-NumFromPlaceholder "C27" 1 "F27"
Set-NumFromPlaceholder "D27" 2 "G27"
Set-NumFromPlaceholder "E27" -50 "H27"
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 6
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = -40
$ws.Range("L27").Value = -53.846153846153
Set-NumFromPlaceholder "C28" 1 "F28"
$ws.Range("F28").Value = 5
$ws.Range("I28").Value = 19
$ws.Range("K28").Value = 137.5
$ws.Range("L28").Value = 72.727272727272
